$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (qol_self_image_r / qol_badadl_r=0.557) is removed: qol_badadl_r now
# occupies row 10 in its own right and the table shrinks from 13 to 12 items.
$ws.Rows.Item(14).Delete()

# Row 2: qol_unpredict_r
$ws.Range("A2").Value = "qol_unpredict_r"
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = "'0.43"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = 0.851

# Row 3: qol_effect_medication_r
$ws.Range("A3").Value = "qol_effect_medication_r"
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "'"
$ws.Range("D3").Value = "'0.9"
$ws.Range("E3").Value = "'"
$ws.Range("F3").Value = 0.249

# Row 4: qol_taking_medication_r
$ws.Range("A4").Value = "qol_taking_medication_r"
$ws.Range("B4").Value = "'"
$ws.Range("C4").Value = "'"
$ws.Range("D4").Value = "'0.79"
$ws.Range("E4").Value = "'"
$ws.Range("F4").Value = 0.363

# Row 5: qol_family_contact_r
$ws.Range("A5").Value = "qol_family_contact_r"
$ws.Range("B5").Value = "'0.79"
$ws.Range("C5").Value = "'"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "'"
$ws.Range("F5").Value = 0.373

# Row 6: qol_friends_contact_r
$ws.Range("A6").Value = "qol_friends_contact_r"
$ws.Range("B6").Value = "'0.96"
$ws.Range("C6").Value = "'"
$ws.Range("D6").Value = "'"
$ws.Range("E6").Value = "'"
$ws.Range("F6").Value = 0.155

# Row 7: qol_affectp_contact_r
$ws.Range("A7").Value = "qol_affectp_contact_r"
$ws.Range("B7").Value = "'0.43"
$ws.Range("C7").Value = "'"
$ws.Range("D7").Value = "'"
$ws.Range("E7").Value = "'"
$ws.Range("F7").Value = 0.816

# Row 8: qol_support_yes_r
$ws.Range("A8").Value = "qol_support_yes_r"
$ws.Range("B8").Value = "'0.47"
$ws.Range("C8").Value = "'"
$ws.Range("D8").Value = "'"
$ws.Range("E8").Value = "'"
$ws.Range("F8").Value = 0.666

# Row 9: qol_badmob_r
$ws.Range("A9").Value = "qol_badmob_r"
$ws.Range("B9").Value = "'"
$ws.Range("C9").Value = "'"
$ws.Range("D9").Value = "'"
$ws.Range("E9").Value = "'0.62"
$ws.Range("F9").Value = 0.613

# Row 10: qol_badadl_r
$ws.Range("A10").Value = "qol_badadl_r"
$ws.Range("B10").Value = "'"
$ws.Range("C10").Value = "'"
$ws.Range("D10").Value = "'"
$ws.Range("E10").Value = "'0.88"
$ws.Range("F10").Value = 0.252

# Row 11: qol_sum_ident_r
$ws.Range("A11").Value = "qol_sum_ident_r"
$ws.Range("B11").Value = "'"
$ws.Range("C11").Value = "'0.7"
$ws.Range("D11").Value = "'"
$ws.Range("E11").Value = "'"
$ws.Range("F11").Value = 0.416

# Row 12: qol_thoughts_future_r
$ws.Range("A12").Value = "qol_thoughts_future_r"
$ws.Range("B12").Value = "'"
$ws.Range("C12").Value = "'0.91"
$ws.Range("D12").Value = "'"
$ws.Range("E12").Value = "'"
$ws.Range("F12").Value = 0.272

# Row 13: qol_self_image_r
$ws.Range("A13").Value = "qol_self_image_r"
$ws.Range("B13").Value = "'"
$ws.Range("C13").Value = "'0.53"
$ws.Range("D13").Value = "'"
$ws.Range("E13").Value = "'"
$ws.Range("F13").Value = 0.592
